$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.820.91"
$ws.Range("E2").Value = "  -1.66%  "

$ws.Range("D3").Value = "3.489.07"
$ws.Range("E3").Value = "  -1.42%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.67"
$ws.Range("E5").Value = "  -1.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.36"
$ws.Range("E6").Value = "  -2.72%  "

$ws.Range("D7").Value = "3.487.43"
$ws.Range("E7").Value = "  -1.39%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.478"
$ws.Range("E9").Value = "  -1.93%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.142"
$ws.Range("E10").Value = "  -0.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.96"
$ws.Range("E11").Value = "  +5.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.423"
$ws.Range("E12").Value = "  -2.81%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000213"
$ws.Range("E13").Value = "  -2.49%  "

$ws.Range("D14").Value = "4.078.03"
$ws.Range("E14").Value = "  -1.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.32"
$ws.Range("E15").Value = "  -4.59%  "

$ws.Range("D16").Value = "3.485.96"
$ws.Range("E16").Value = "  -1.36%  "

$ws.Range("D17").Value = "66.858.98"
$ws.Range("E17").Value = "  -1.62%  "

$ws.Range("E18").Value = "  +0.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.39"
$ws.Range("E19").Value = "  -3.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.28"
$ws.Range("E20").Value = "  +3.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.30"
$ws.Range("E21").Value = "  -2.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "434.14"
$ws.Range("E22").Value = "  -3.79%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.608"
$ws.Range("E23").Value = "  -4.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.70"
$ws.Range("E24").Value = "  +2.00%  "

$ws.Range("E25").Value = "  +0.09%  "

$ws.Range("D26").Value = "3.624.25"
$ws.Range("E26").Value = "  -1.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000119"
$ws.Range("E27").Value = "  -8.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.82"
$ws.Range("E28").Value = "  -4.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.20"
$ws.Range("E29").Value = "  -9.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.48"
$ws.Range("E30").Value = "  -1.99%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.61"
$ws.Range("E31").Value = "  -5.03%  "

$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("E33").Value = "  -3.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.42"
$ws.Range("E34").Value = "  -1.95%  "

$ws.Range("D35").Value = "3.480.28"
$ws.Range("E35").Value = "  -1.51%  "

$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.95"
$ws.Range("E36").Value = "  -4.92%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.81"
$ws.Range("E37").Value = "  -4.84%  "

$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.01"
$ws.Range("E38").Value = "  -1.18%  "

$ws.Range("B39").Value = "USDe"
$ws.Range("C39").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.15%  "

$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0892"
$ws.Range("E41").Value = "  -1.91%  "

$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "170.67"
$ws.Range("E42").Value = "  -3.72%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.43"
$ws.Range("E43").Value = "  -2.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.08"
$ws.Range("E44").Value = "  -10.79%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.897"
$ws.Range("E45").Value = "  +0.85%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.15"
$ws.Range("E46").Value = "  -6.80%  "

$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.85"
$ws.Range("E47").Value = "  -0.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.24"
$ws.Range("E48").Value = "  -6.76%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.46"
$ws.Range("E49").Value = "  -3.27%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.42"
$ws.Range("E50").Value = "  -6.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.970"
$ws.Range("E51").Value = "  -3.53%  "
